$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (B3/C3) already carries the "leftover hyperlink" cell format
# (plain text style that used to back a hyperlink). Copy that same
# formatting onto row 2's B2/C2 cells so all four cells share one style,
# using a formats-only paste so the existing style is reused rather than
# a new one minted.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Swap the team names that sit in column C between rows 2 and 3.
$ws.Range("C2").Value = "Toronto Maple Leafs"
$ws.Range("C3").Value = "Montreal Canadiens"

# The two hyperlinks (on B2 and C2) are removed, leaving the plain
# values/formatting behind.
$ws.Hyperlinks.Delete() | Out-Null

# Move the sheet's active selection from C10 to C8.
$ws.Range("C8").Select() | Out-Null
